# "Generate Report for Handoff" -- update the status row for the
# eba56d22-c6f5-4b2c-b9c7-83d400b3b666 file (row 6 in every sheet) to
# reflect that it is now "Ready for handoff".

$wb = $excel.ActiveWorkbook

$statusReadyForHandoff = "Ready for handoff"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test4/blob/5cc56fa0a4458d51a8e52ab65af7773ed591c2dc/e2e/eba56d22-c6f5-4b2c-b9c7-83d400b3b666.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test4/blob/3668b2b0c4595bc910ca8df8459fd51e02a5ffdd/e2e/eba56d22-c6f5-4b2c-b9c7-83d400b3b666.md."

# --- Overview sheet --------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E6").Value = $statusReadyForHandoff
$wsOverview.Range("F6").Value = $statusReadyForHandoff
$wsOverview.Range("G6").Value = "2017-02-28 08:07:13"

# --- zh-cn sheet -------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C6").Value = $statusReadyForHandoff
$wsZhCn.Range("E6").Value = "ht"
$wsZhCn.Range("H6").Value = "2017-02-28 08:06:56"
$wsZhCn.Range("R6").Value = $errorDetail
$wsZhCn.Columns.Item(18).ColumnWidth = 40

# --- de-de sheet -------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C6").Value = $statusReadyForHandoff
$wsDeDe.Range("E6").Value = "ht"
$wsDeDe.Range("H6").Value = "2017-02-28 08:07:13"
$wsDeDe.Range("R6").Value = $errorDetail
$wsDeDe.Columns.Item(18).ColumnWidth = 40
